$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.357.66"
$ws.Range("E2").Value = "  +0.15%  "
$ws.Range("D3").Value = "'1.876.50"
$ws.Range("E3").Value = "  +0.38%  "
$ws.Range("D4").Value = "'0.9999"
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").Value = "'0.7140"
$ws.Range("E5").Value = "  -1.56%  "
$ws.Range("D6").Value = "'242.09"
$ws.Range("E6").Value = "  +0.36%  "
$ws.Range("D8").Value = "'0.3115"
$ws.Range("E8").Value = "  +0.71%  "
$ws.Range("D9").Value = "'0.07711"
$ws.Range("E9").Value = "  -2.12%  "
$ws.Range("D10").Value = "'25.18"
$ws.Range("E10").Value = "  -0.29%  "
$ws.Range("D11").Value = "'0.08378"
$ws.Range("E11").Value = "  +1.35%  "
$ws.Range("D12").Value = "'1.888.17"
$ws.Range("E12").Value = "  +0.17%  "
$ws.Range("E13").Value = "  -0.11%  "
$ws.Range("D14").Value = "'0.7145"
$ws.Range("E14").Value = "  -1.12%  "
$ws.Range("D15").Value = "'91.65"
$ws.Range("E15").Value = "  +1.07%  "
$ws.Range("D16").Value = "'0.000008400"
$ws.Range("E16").Value = "  +7.36%  "
$ws.Range("D17").Value = "'29.371.54"
$ws.Range("E17").Value = "  +0.15%  "
$ws.Range("D18").Value = "'5.969"
$ws.Range("E18").Value = "  +2.01%  "
$ws.Range("D19").Value = "'243.17"
$ws.Range("E19").Value = "  -0.37%  "
$ws.Range("D20").Value = "'2.133.84"
$ws.Range("E20").Value = "  +1.22%  "
$ws.Range("D21").Value = "'13.21"
$ws.Range("E21").Value = "  -0.09%  "
$ws.Range("D22").Value = "'0.9994"
$ws.Range("E22").Value = "  -0.12%  "
$ws.Range("D23").Value = "'7.883"
$ws.Range("E23").Value = "  -1.08%  "
$ws.Range("E24").Value = "  -0.06%  "
$ws.Range("D25").Value = "'0.1619"
$ws.Range("E25").Value = "  +1.04%  "
$ws.Range("D26").Value = "'164.12"
$ws.Range("E26").Value = "  +0.80%  "
$ws.Range("D27").Value = "'9.013"
$ws.Range("E27").Value = "  +0.58%  "
$ws.Range("D28").Value = "'18.57"
$ws.Range("E28").Value = "  +1.68%  "
$ws.Range("E29").Value = "  +0.76%  "
$ws.Range("D30").Value = "'4.409"
$ws.Range("E30").Value = "  +0.47%  "
$ws.Range("D31").Value = "'4.326"
$ws.Range("E31").Value = "  +5.05%  "
$ws.Range("D32").Value = "'1.289"
$ws.Range("E32").Value = "  -4.55%  "
$ws.Range("D33").Value = "'0.05245"
$ws.Range("E33").Value = "  +0.62%  "
$ws.Range("D34").Value = "'1.924"
$ws.Range("E34").Value = "  -0.52%  "
$ws.Range("D35").Value = "'0.7598"
$ws.Range("E35").Value = "  +4.43%  "
$ws.Range("D36").Value = "'1.174"
$ws.Range("E36").Value = "  -1.09%  "
$ws.Range("D37").Value = "'2.677"
$ws.Range("E37").Value = "  -0.25%  "
$ws.Range("D38").Value = "'0.01863"
$ws.Range("E38").Value = "  +0.22%  "
$ws.Range("D39").Value = "'2.721"
$ws.Range("E39").Value = "  +0.77%  "
$ws.Range("D40").Value = "'1.160.31"
$ws.Range("E40").Value = "  -0.84%  "
$ws.Range("D41").Value = "'6.345"
$ws.Range("E41").Value = "  +3.84%  "
$ws.Range("D42").Value = "'73.42"
$ws.Range("E42").Value = "  +0.90%  "
$ws.Range("D43").Value = "'0.8897"
$ws.Range("E43").Value = "  -1.59%  "
$ws.Range("D44").Value = "'104.70"
$ws.Range("E44").Value = "  +2.91%  "
$ws.Range("D45").Value = "'0.9996"
$ws.Range("E45").Value = "  -0.19%  "
$ws.Range("D46").Value = "'2.032.57"
$ws.Range("E46").Value = "  +1.26%  "
$ws.Range("D47").Value = "'0.5201"
$ws.Range("E47").Value = "  -1.59%  "
$ws.Range("D48").Value = "'1.795"
$ws.Range("E48").Value = "  +0.72%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "'9.411"
$ws.Range("E49").Value = "  +1.47%  "
$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").Value = "'0.00000000120"
$ws.Range("E50").Value = "  -0.47%  "
$ws.Range("D51").Value = "'0.4306"
$ws.Range("E51").Value = "  +0.53%  "
